# "Changes of 27th April 2022"
#
# The MDSi test-result sheet gets re-run with a new batch of Job# values.
# Rows 2-4 (column B, "Job#") are refreshed with the three newest job
# numbers from that day's run:
#   B2: 32341649 -> 32372249
#   B3: 32341650 -> 32372251
#   B4: 32341651 -> 32372252
#
# These Job# values are digit-only but must stay TEXT cells (as they were
# before the edit - shared-string typed, not numeric), matching how the
# column already stored them. Assigning a numeric-looking string straight
# to .Value would auto-coerce it to a Number, so the cells are briefly
# switched to a text number format, written, then restored to the sheet's
# normal (General) formatting so no visible formatting change remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered (not a hashtable) so B2, B3, B4 are written in exactly this
# sequence.
$jobNumberUpdates = @(
    , @("B2", "32372249")
    , @("B3", "32372251")
    , @("B4", "32372252")
)

foreach ($update in $jobNumberUpdates) {
    $cellRef = $update[0]
    $newJobNumber = $update[1]

    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newJobNumber
    $cell.Style = "Normal"
}
